$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "316.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.61%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.62%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.158"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.11%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07968"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.01%"

$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.482"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.66%"

$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.922"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.92%"

$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.980"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "6.61%"

$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9408"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.32%"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1275"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "6.97%"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1940"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.87%"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09055"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.69%"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03421"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.00%"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09527"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.87%"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001371"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.96%"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006072"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "5.04%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.420"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.80%"

$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.464"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.86%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3512"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.06%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.554"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "24.41%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1305"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.54%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2300"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-11.41%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04359"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.30%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001198"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-4.45%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004409"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.44%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001325"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-2.83%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003971"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.93%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02386"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.38%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05176"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.34%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007426"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.55%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1396"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.38%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008358"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-8.08%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002102"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "7.46%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008738"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006471"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.43%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.93%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002852"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-13.26%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001681"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "67.37%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002090"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.93%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001990"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.93%"
